$wb = $excel.ActiveWorkbook
$wsF = $wb.Worksheets.Item("F")
$wsOutput = $wb.Worksheets.Item("OUTPUT")
$wsSkills = $wb.Worksheets.Item("SkillFormulas")

# --- Content edits on sheet F (character stats updated to a new morph) ---

# Skill name: Computer Programming -> Fabber Programming
$wsF.Range("A25").Value = "Fabber Programming"

# Morph description line 2 (appearance)
$wsF.Range("J2").ClearFormats()
$wsF.Range("J2").Value = "bald, 158 cm, 54 kg"

# Morph name
$wsF.Range("A4").Value = "Ocon SFF Lo-Grav"

# Morph traits
$wsF.Range("A7").Value = "DR 3, Mesh Inserts, Cyberbrain, Mnemonic Augmentation, Cortical Stack, Access Jacks"

# Morph stat bonuses (row 4: ST, DX, IQ, WILL, PER, HT)
$wsF.Range("C4").Value = 10
$wsF.Range("D4").Value = 11
$wsF.Range("E4").Value = 10
$wsF.Range("F4").Value = 10
$wsF.Range("G4").Value = 10
$wsF.Range("H4").Value = 10

# --- Restore view/selection state ---
$wsOutput.Range("A18").Select()
$wsSkills.Range("G5").Select()
$wsF.Select()
$wsF.Range("B12").Select()
